# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the "Estado de Cuenta" worker/period detail table (rows 16-37):
#   - JAIME LUIS ORTEGA GARCIA keeps periods 2003,2005-2012 in rows 16-24
#     (value for 2003 now uses the 24578 amount that used to sit on 2106)
#   - Rows 25-37 now interleave CARLOS CABALLERO MONTES / JAIME LUIS ORTEGA
#     GARCIA for periods 2101-2106 (new "parte 1" rows for Carlos)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2003"; Mora = 24578;  Salario = 877803 },
    @{ Row = 17; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2005"; Mora = 35112;  Salario = 877803 },
    @{ Row = 18; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2006"; Mora = 35112;  Salario = 877803 },
    @{ Row = 19; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2007"; Mora = 35112;  Salario = 877803 },
    @{ Row = 20; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2008"; Mora = 35112;  Salario = 877803 },
    @{ Row = 21; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2009"; Mora = 35112;  Salario = 877803 },
    @{ Row = 22; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2010"; Mora = 35112;  Salario = 877803 },
    @{ Row = 23; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2011"; Mora = 35112;  Salario = 877803 },
    @{ Row = 24; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2012"; Mora = 35112;  Salario = 877803 },
    @{ Row = 25; Doc = "1047408531"; Nombre = "CARLOS CABALLERO MONTES";  Periodo = "2012"; Mora = 36000;  Salario = 900000 },
    @{ Row = 26; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2101"; Mora = 35112;  Salario = 877803 },
    @{ Row = 27; Doc = "1047408531"; Nombre = "CARLOS CABALLERO MONTES";  Periodo = "2101"; Mora = 36000;  Salario = 900000 },
    @{ Row = 28; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2102"; Mora = 35112;  Salario = 877803 },
    @{ Row = 29; Doc = "1047408531"; Nombre = "CARLOS CABALLERO MONTES";  Periodo = "2102"; Mora = 36000;  Salario = 900000 },
    @{ Row = 30; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2103"; Mora = 35112;  Salario = 877803 },
    @{ Row = 31; Doc = "1047408531"; Nombre = "CARLOS CABALLERO MONTES";  Periodo = "2103"; Mora = 36000;  Salario = 900000 },
    @{ Row = 32; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2104"; Mora = 35112;  Salario = 877803 },
    @{ Row = 33; Doc = "1047408531"; Nombre = "CARLOS CABALLERO MONTES";  Periodo = "2104"; Mora = 36000;  Salario = 900000 },
    @{ Row = 34; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2105"; Mora = 35112;  Salario = 877803 },
    @{ Row = 35; Doc = "1047408531"; Nombre = "CARLOS CABALLERO MONTES";  Periodo = "2105"; Mora = 36000;  Salario = 900000 },
    @{ Row = 36; Doc = "1143390056"; Nombre = "JAIME LUIS ORTEGA GARCIA"; Periodo = "2106"; Mora = 24578;  Salario = 877803 },
    @{ Row = 37; Doc = "1047408531"; Nombre = "CARLOS CABALLERO MONTES";  Periodo = "2106"; Mora = 25200;  Salario = 900000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = "CC"
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Nombre
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
